$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A2" = -0.0
    "B2" = -0.07498804478836586
    "C2" = -0.0
    "D2" = 0.2123108505157703
    "E2" = 0.005516432524130028
    "G2" = 0.0
    "I2" = -0.0
    "J2" = -0.0
    "K2" = 0.0210650617049058
    "L2" = -0.0
    "M2" = 0.2062088150060369
    "N2" = -0.007387703804133994
    "R2" = -0.0
    "S2" = 0.0
    "T2" = -0.09179633776493132
    "V2" = 0.01626758625129003
    "W2" = -0.03827681932504355
    "Y2" = -0.0
    "Z2" = -0.0
    "AB2" = 0.0
    "AC2" = -0.05304819134633224
    "AD2" = 0.0
    "AE2" = -0.01596718197831162
    "AF2" = 0.001357506954384982
    "AG2" = -0.0
    "AH2" = -0.0
    "AI2" = -0.0
    "AJ2" = 0.0
    "AK2" = -0.0
    "AL2" = -0.03337734227611048
    "AM2" = 0.0
    "AN2" = 0.02710235950979624
    "AO2" = 0.06807114299100285
    "AQ2" = 0.0
    "AR2" = -0.0
    "AT2" = 0.0
    "AU2" = -0.1493238582359526
    "AW2" = 0.07989574898010814
    "AX2" = 0.004724416540976223
    "AY2" = -0.0
    "BC2" = -0.0
    "BD2" = -0.01731206187656933
    "BF2" = 0.08423891762249715
    "BG2" = 0.03181235655624312
    "BJ2" = -0.0
    "BL2" = 0.0
    "BM2" = 0.03624722305829137
    "BO2" = -0.03913293478256753
    "BP2" = -0.08637971709137597
    "BU2" = 0.0
    "BV2" = -0.04500391369559421
    "BW2" = 0.0
    "BX2" = 0.01500758071599884
    "BY2" = -0.02086773070508119
    "BZ2" = -0.0
    "CB2" = 0.0
    "CD2" = -0.0
    "CE2" = 0.03220467152992416
    "CG2" = -0.03364726260395302
    "CH2" = 0.01645725400885414
    "CJ2" = -0.0
    "CM2" = -0.0
    "CN2" = -0.01104918755697818
    "CP2" = 0.02074356515685865
    "CQ2" = 0.03582172844414858
    "CT2" = 0.0
    "CU2" = -0.0
    "CV2" = -0.0
    "CW2" = 0.04603546316646197
    "CY2" = -0.03621017246313984
    "CZ2" = 0.009757611166591691
    "DD2" = -0.0
    "DE2" = -0.0
    "DF2" = 0.02952187127600957
    "DH2" = 0.02516906389234852
    "DI2" = 0.03783215140409813
    "DJ2" = 0.0
    "DK2" = -0.0
    "DL2" = -0.0
    "DN2" = 0.0
    "DO2" = -0.01844302020456456
    "DQ2" = 0.03409103206182655
    "DR2" = -0.01622991875651391
    "DS2" = -0.0
    "DW2" = 0.0
    "DX2" = -0.0547282908726129
    "DY2" = -0.0
    "DZ2" = -0.01156697219937911
    "EA2" = -0.02323154868926323
    "EB2" = 0.0
    "EF2" = -0.0
    "EG2" = 0.04188967657062639
    "EI2" = 0.07034029962704759
    "EJ2" = -0.02720351628130804
    "EO2" = 0.0
    "EP2" = 0.04467283160219787
    "EQ2" = 0.0
    "ER2" = -0.0332712947545084
    "ES2" = 0.03374213308970007
    "ET2" = 0.0
    "EU2" = -0.0
    "EV2" = 0.0
    "EX2" = 0.0
    "EY2" = 0.04336569382910841
    "FA2" = -0.02930277627662451
    "FB2" = 0.01670699430381313
    "FD2" = -0.0
    "FG2" = -0.0
    "FH2" = 0.003066748497303185
    "FI2" = 0.0
    "FJ2" = -0.006813118420519349
    "FK2" = -0.006099503072460796
    "FL2" = -0.0
    "FN2" = -0.0
    "FP2" = -0.0
    "FQ2" = -0.01240537822611858
    "FR2" = -0.0
    "FS2" = -0.02042875914825069
    "FT2" = 0.007039864533620131
    "FV2" = -0.0
    "FW2" = -0.0
    "FY2" = 0.0
    "FZ2" = -0.03029568608595597
    "GB2" = 0.03483274327951914
    "GD2" = 0.0
    "GE2" = -0.0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
